$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Nos" — node coordinates table grows from 3 rows (A2:B4) to 9
# rows (A2:B10); D2 keeps its COUNT() formula which will recalc to 9.
# ---------------------------------------------------------------------
$wsNos = $wb.Worksheets.Item("Nos")

$wsNos.Cells.Item(2,1).Value = 0
$wsNos.Cells.Item(2,2).Value = 0

$wsNos.Cells.Item(3,1).Value = 1
$wsNos.Cells.Item(3,2).Value = 0

$wsNos.Cells.Item(4,1).Value = 2
$wsNos.Cells.Item(4,2).Value = 0

# Rows 5-10 are brand new — copy the formatting that row 5 already
# carries (style used for the bottom of the table) down before writing
# the values so the new rows pick up matching look & feel.
$wsNos.Range("A5:B5").Copy()
$wsNos.Range("A5:B10").PasteSpecial(-4122)

$wsNos.Cells.Item(5,1).Value = 3
$wsNos.Cells.Item(5,2).Value = 0

$wsNos.Cells.Item(6,1).Value = 4
$wsNos.Cells.Item(6,2).Value = 0

$wsNos.Cells.Item(7,1).Value = 3
$wsNos.Cells.Item(7,2).Value = 1

$wsNos.Cells.Item(8,1).Value = 2
$wsNos.Cells.Item(8,2).Value = 2

$wsNos.Cells.Item(9,1).Value = 1
$wsNos.Cells.Item(9,2).Value = 1

$wsNos.Cells.Item(10,1).Value = 2
$wsNos.Cells.Item(10,2).Value = 1

$wsNos.Range("E12").Select()

# ---------------------------------------------------------------------
# Sheet "Incidencia" — member list grows from 3 rows (A2:D4) to 16 rows
# (A2:D17). C column stays 210e9 throughout; D gets a real formula
# (=0.15*0.15) instead of the hard-coded 2E-4 literal, shared across
# D3:D17.
# ---------------------------------------------------------------------
$wsInc = $wb.Worksheets.Item("Incidencia")

$wsInc.Cells.Item(2,1).Value = 1
$wsInc.Cells.Item(2,2).Value = 2
$wsInc.Cells.Item(2,3).Value = 210000000000
$wsInc.Range("D2").Formula = "=0.15*0.15"

$wsInc.Cells.Item(3,1).Value = 2
$wsInc.Cells.Item(3,2).Value = 3
$wsInc.Cells.Item(3,3).Value = 210000000000

$wsInc.Cells.Item(4,1).Value = 3
$wsInc.Cells.Item(4,2).Value = 4
$wsInc.Cells.Item(4,3).Value = 210000000000

# Copy row 4's formatting down across the new rows before filling them.
$wsInc.Range("A4:D4").Copy()
$wsInc.Range("A5:D17").PasteSpecial(-4122)

$wsInc.Cells.Item(5,1).Value = 4
$wsInc.Cells.Item(5,2).Value = 5
$wsInc.Cells.Item(5,3).Value = 210000000000

$wsInc.Cells.Item(6,1).Value = 5
$wsInc.Cells.Item(6,2).Value = 6
$wsInc.Cells.Item(6,3).Value = 210000000000

$wsInc.Cells.Item(7,1).Value = 6
$wsInc.Cells.Item(7,2).Value = 7
$wsInc.Cells.Item(7,3).Value = 210000000000

$wsInc.Cells.Item(8,1).Value = 7
$wsInc.Cells.Item(8,2).Value = 8
$wsInc.Cells.Item(8,3).Value = 210000000000

$wsInc.Cells.Item(9,1).Value = 8
$wsInc.Cells.Item(9,2).Value = 1
$wsInc.Cells.Item(9,3).Value = 210000000000

$wsInc.Cells.Item(10,1).Value = 2
$wsInc.Cells.Item(10,2).Value = 8
$wsInc.Cells.Item(10,3).Value = 210000000000

$wsInc.Cells.Item(11,1).Value = 2
$wsInc.Cells.Item(11,2).Value = 9
$wsInc.Cells.Item(11,3).Value = 210000000000

$wsInc.Cells.Item(12,1).Value = 3
$wsInc.Cells.Item(12,2).Value = 9
$wsInc.Cells.Item(12,3).Value = 210000000000

$wsInc.Cells.Item(13,1).Value = 9
$wsInc.Cells.Item(13,2).Value = 4
$wsInc.Cells.Item(13,3).Value = 210000000000

$wsInc.Cells.Item(14,1).Value = 4
$wsInc.Cells.Item(14,2).Value = 6
$wsInc.Cells.Item(14,3).Value = 210000000000

$wsInc.Cells.Item(15,1).Value = 9
$wsInc.Cells.Item(15,2).Value = 6
$wsInc.Cells.Item(15,3).Value = 210000000000

$wsInc.Cells.Item(16,1).Value = 8
$wsInc.Cells.Item(16,2).Value = 9
$wsInc.Cells.Item(16,3).Value = 210000000000

$wsInc.Cells.Item(17,1).Value = 9
$wsInc.Cells.Item(17,2).Value = 7
$wsInc.Cells.Item(17,3).Value = 210000000000

# D3:D17 share the same "=0.15*0.15" formula (written as one shot so
# Excel records it as a single shared-formula group, matching the
# original authoring).
$wsInc.Range("D3:D17").Formula = "=0.15*0.15"

$wsInc.Range("A1:A1048576").Select()

# ---------------------------------------------------------------------
# Sheet "Carregamento" — loads table: values change in rows 2-3 and a
# new row 4 is added.
# ---------------------------------------------------------------------
$wsCar = $wb.Worksheets.Item("Carregamento")

$wsCar.Cells.Item(2,1).Value = 7
$wsCar.Cells.Item(2,2).Value = 1
$wsCar.Cells.Item(2,3).Value = 2000

$wsCar.Cells.Item(3,1).Value = 7
$wsCar.Cells.Item(3,2).Value = 2
$wsCar.Cells.Item(3,3).Value = -10000

$wsCar.Cells.Item(4,1).Value = 8
$wsCar.Cells.Item(4,2).Value = 1
$wsCar.Cells.Item(4,3).Value = 2000

$wsCar.Range("A4:C4").Select()

# ---------------------------------------------------------------------
# Sheet "Restricao" — constraints table: row 3 values swap, and rows
# 5-8 gain new node/direction pairs.
# ---------------------------------------------------------------------
$wsRes = $wb.Worksheets.Item("Restricao")

$wsRes.Cells.Item(2,1).Value = 1
$wsRes.Cells.Item(2,2).Value = 1

$wsRes.Cells.Item(3,1).Value = 1
$wsRes.Cells.Item(3,2).Value = 2

$wsRes.Cells.Item(4,1).Value = 2
$wsRes.Cells.Item(4,2).Value = 1

# A5:B5 are brand new cells — borrow the formatting already sitting on
# A6:B6 (same look as the rest of the table) before writing into them.
$wsRes.Range("A6:B6").Copy()
$wsRes.Range("A5:B5").PasteSpecial(-4122)

$wsRes.Cells.Item(5,1).Value = 3
$wsRes.Cells.Item(5,2).Value = 2

$wsRes.Cells.Item(6,1).Value = 4
$wsRes.Cells.Item(6,2).Value = 2

$wsRes.Cells.Item(7,1).Value = 5
$wsRes.Cells.Item(7,2).Value = 1

$wsRes.Cells.Item(8,1).Value = 5
$wsRes.Cells.Item(8,2).Value = 2

$wsRes.Range("C12").Select()

# ---------------------------------------------------------------------
# Make "Restricao" the active sheet/tab (it was tabSelected in both the
# before and after states) and recalc everything.
# ---------------------------------------------------------------------
$wsRes.Activate()
$excel.Calculate()
